$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    [double]"3.851225878354915e-10",
    [double]"7.020160124933402e-12",
    [double]"1.170456251936844e-13",
    [double]"1.562852079366264e-15",
    [double]"2.848825853191154e-17",
    [double]"6.343225018861762e-19",
    [double]"2.108688116291305e-20",
    [double]"3.843796416384833e-22",
    [double]"6.408679384322087e-24",
    [double]"134962.5154844009",
    [double]"1757832.594709143",
    [double]"768538.6203265248",
    [double]"974563.6817176824",
    [double]"452613.5335914007",
    [double]"26198.84397427095",
    [double]"16269835.90682306",
    [double]"15702821.79366891",
    [double]"8331903.210997708",
    [double]"7629416.52609037",
    [double]"4304435.526403553",
    [double]"1612123.061485417",
    [double]"1124896.296864709",
    [double]"2849689.676257911",
    [double]"788393.6281928831",
    [double]"4987356.099866037",
    [double]"10834204.15497985",
    [double]"348766.2390630613",
    [double]"6357.442855154855",
    [double]"115.8858717665375",
    [double]"1.932140303161275",
    [double]"0.02579890948872387",
    [double]"243074.544019867",
    [double]"484519.0046972869",
    [double]"8831.995587860069",
    [double]"160.9929544718569",
    [double]"2.684201025701763"
)

$row = 2
foreach ($val in $values) {
    $ws.Cells.Item($row, 3).Value = $val
    $row = $row + 1
}

$wb.Save()
